$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 20,10
$data[0,0] = -18.58334357994523
$data[0,1] = -1.060356010718436
$data[0,2] = -18.58334357994523
$data[0,3] = -18.58334357994523
$data[0,4] = -18.58334357994523
$data[0,5] = -18.58334357994523
$data[0,6] = -18.58334357994523
$data[0,7] = -18.58334357994523
$data[0,8] = -18.58334357994523
$data[0,9] = -18.58334357994523
$data[1,0] = -18.58334357994523
$data[1,1] = -18.58334357994523
$data[1,2] = -18.58334357994523
$data[1,3] = -18.58334357994523
$data[1,4] = -18.58334357994523
$data[1,5] = -18.58334357994523
$data[1,6] = -18.58334357994523
$data[1,7] = 0.2010899827321729
$data[1,8] = -18.58334357994523
$data[1,9] = -18.58334357994523
$data[2,0] = -18.58334357994523
$data[2,1] = -1.006231969612305
$data[2,2] = -0.5820398527712781
$data[2,3] = -18.58334357994523
$data[2,4] = 3.859419998971665
$data[2,5] = -18.58334357994523
$data[2,6] = 1.813166088123149
$data[2,7] = -18.58334357994523
$data[2,8] = 2.559723815610687
$data[2,9] = -18.58334357994523
$data[3,0] = -18.58334357994523
$data[3,1] = 0.3063653511920159
$data[3,2] = -18.58334357994523
$data[3,3] = -18.58334357994523
$data[3,4] = -18.58334357994523
$data[3,5] = 3.463831266743902
$data[3,6] = -18.58334357994523
$data[3,7] = -18.58334357994523
$data[3,8] = -18.58334357994523
$data[3,9] = -18.58334357994523
$data[4,0] = -18.58334357994523
$data[4,1] = -18.58334357994523
$data[4,2] = -18.58334357994523
$data[4,3] = -18.58334357994523
$data[4,4] = -18.58334357994523
$data[4,5] = -18.58334357994523
$data[4,6] = -18.58334357994523
$data[4,7] = -18.58334357994523
$data[4,8] = -18.58334357994523
$data[4,9] = -18.58334357994523
$data[5,0] = 3.158141922600064
$data[5,1] = -18.58334357994523
$data[5,2] = -18.58334357994523
$data[5,3] = -18.58334357994523
$data[5,4] = -18.58334357994523
$data[5,5] = -18.58334357994523
$data[5,6] = -18.58334357994523
$data[5,7] = -18.58334357994523
$data[5,8] = -18.58334357994523
$data[5,9] = -18.58334357994523
$data[6,0] = -18.58334357994523
$data[6,1] = -18.58334357994523
$data[6,2] = -18.58334357994523
$data[6,3] = -18.58334357994523
$data[6,4] = -18.58334357994523
$data[6,5] = -18.58334357994523
$data[6,6] = -18.58334357994523
$data[6,7] = -18.58334357994523
$data[6,8] = -18.58334357994523
$data[6,9] = -18.58334357994523
$data[7,0] = 3.468995277189252
$data[7,1] = -18.58334357994523
$data[7,2] = -18.58334357994523
$data[7,3] = -18.58334357994523
$data[7,4] = -18.58334357994523
$data[7,5] = -18.58334357994523
$data[7,6] = -18.58334357994523
$data[7,7] = -18.58334357994523
$data[7,8] = -18.58334357994523
$data[7,9] = -18.58334357994523
$data[8,0] = -18.58334357994523
$data[8,1] = -18.58334357994523
$data[8,2] = -18.58334357994523
$data[8,3] = -18.58334357994523
$data[8,4] = -18.58334357994523
$data[8,5] = -18.58334357994523
$data[8,6] = -18.58334357994523
$data[8,7] = 0.4047558029249811
$data[8,8] = -18.58334357994523
$data[8,9] = 2.212322761611524
$data[9,0] = -18.58334357994523
$data[9,1] = -18.58334357994523
$data[9,2] = -18.58334357994523
$data[9,3] = -18.58334357994523
$data[9,4] = -18.58334357994523
$data[9,5] = 1.626948445355757
$data[9,6] = -18.58334357994523
$data[9,7] = -18.58334357994523
$data[9,8] = -18.58334357994523
$data[9,9] = 1.585569655942634
$data[10,0] = -18.58334357994523
$data[10,1] = -18.58334357994523
$data[10,2] = -18.58334357994523
$data[10,3] = -18.58334357994523
$data[10,4] = -18.58334357994523
$data[10,5] = -18.58334357994523
$data[10,6] = -18.58334357994523
$data[10,7] = -18.58334357994523
$data[10,8] = -18.58334357994523
$data[10,9] = -18.58334357994523
$data[11,0] = -18.58334357994523
$data[11,1] = -18.58334357994523
$data[11,2] = -18.58334357994523
$data[11,3] = -18.58334357994523
$data[11,4] = -18.58334357994523
$data[11,5] = -18.58334357994523
$data[11,6] = -18.58334357994523
$data[11,7] = -18.58334357994523
$data[11,8] = 0.953967731749447
$data[11,9] = 2.521749669103986
$data[12,0] = -18.58334357994523
$data[12,1] = -18.58334357994523
$data[12,2] = -0.7145149978806806
$data[12,3] = -18.58334357994523
$data[12,4] = -18.58334357994523
$data[12,5] = -18.58334357994523
$data[12,6] = -18.58334357994523
$data[12,7] = -18.58334357994523
$data[12,8] = -18.58334357994523
$data[12,9] = 1.741752119882642
$data[13,0] = -18.58334357994523
$data[13,1] = -18.58334357994523
$data[13,2] = -1.3679148578159
$data[13,3] = -18.58334357994523
$data[13,4] = -18.58334357994523
$data[13,5] = -18.58334357994523
$data[13,6] = -18.58334357994523
$data[13,7] = -18.58334357994523
$data[13,8] = -18.58334357994523
$data[13,9] = -18.58334357994523
$data[14,0] = -18.58334357994523
$data[14,1] = -18.58334357994523
$data[14,2] = -18.58334357994523
$data[14,3] = -18.58334357994523
$data[14,4] = -18.58334357994523
$data[14,5] = -18.58334357994523
$data[14,6] = -18.58334357994523
$data[14,7] = -18.58334357994523
$data[14,8] = 2.110957338688764
$data[14,9] = -18.58334357994523
$data[15,0] = -18.58334357994523
$data[15,1] = -0.2029121090776649
$data[15,2] = -0.7209180998803044
$data[15,3] = -18.58334357994523
$data[15,4] = -18.58334357994523
$data[15,5] = -18.58334357994523
$data[15,6] = 2.128586327336272
$data[15,7] = 0.0942831653945098
$data[15,8] = 1.988661751514682
$data[15,9] = -18.58334357994523
$data[16,0] = -18.58334357994523
$data[16,1] = -18.58334357994523
$data[16,2] = -18.58334357994523
$data[16,3] = -18.58334357994523
$data[16,4] = -18.58334357994523
$data[16,5] = -18.58334357994523
$data[16,6] = 2.217402314016716
$data[16,7] = -0.5371094819666122
$data[16,8] = 1.955444973430715
$data[16,9] = -18.58334357994523
$data[17,0] = -18.58334357994523
$data[17,1] = -18.58334357994523
$data[17,2] = 3.378377761435282
$data[17,3] = -18.58334357994523
$data[17,4] = -18.58334357994523
$data[17,5] = -18.58334357994523
$data[17,6] = 1.649087887455752
$data[17,7] = 1.370643939491468
$data[17,8] = -18.58334357994523
$data[17,9] = -18.58334357994523
$data[18,0] = -18.58334357994523
$data[18,1] = 3.274452563503753
$data[18,2] = 2.873683031065483
$data[18,3] = -18.58334357994523
$data[18,4] = 2.455618934033639
$data[18,5] = -18.58334357994523
$data[18,6] = 1.105215796224917
$data[18,7] = 3.720714960953494
$data[18,8] = -18.58334357994523
$data[18,9] = 1.712507172661907
$data[19,0] = -18.58334357994523
$data[19,1] = 2.856177397885398
$data[19,2] = -18.58334357994523
$data[19,3] = 4.321924605454617
$data[19,4] = -18.58334357994523
$data[19,5] = 2.555268836901033
$data[19,6] = 1.120838702199291
$data[19,7] = -18.58334357994523
$data[19,8] = -18.58334357994523
$data[19,9] = -18.58334357994523

$ws.Range("B2:K21").Value = $data

Write-Output "done"